$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns: F "Bl" and G "Operating Freq"
$ws.Range("F1").Value = "Bl"
$ws.Range("G1").Value = "Operating Freq"

# Match the header formatting already used by A1:E1 (centered alignment style)
$ws.Range("A1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the active selection to F2
$ws.Range("F2").Select()
